# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (interest/attendance count) figures in column F
# across the 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 90
$ws1.Range("F3").Value = 11945
$ws1.Range("F4").Value = 26
$ws1.Range("F5").Value = 224
$ws1.Range("F8").Value = 11839
$ws1.Range("F11").Value = 106
$ws1.Range("F12").Value = 69
$ws1.Range("F13").Value = 1785
$ws1.Range("F14").Value = 5867
$ws1.Range("F18").Value = 23

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 6

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 90
$ws4.Range("F5").Value = 11945
$ws4.Range("F6").Value = 26
$ws4.Range("F7").Value = 224
$ws4.Range("F8").Value = 6
$ws4.Range("F11").Value = 11839
$ws4.Range("F14").Value = 106
$ws4.Range("F15").Value = 69
$ws4.Range("F16").Value = 1785
$ws4.Range("F18").Value = 5868
$ws4.Range("F22").Value = 23

$wb.Save()
